$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual "missing data" cells in column E (rows 2-23) ---
# These flip which cells are treated as missing (blank) vs populated,
# simulating a different random seed for data removal.
$ws.Cells.Item(2, 5).Value = -7.2          # E2: was blank -> -7.2
$ws.Cells.Item(6, 5).ClearContents()        # E6: was -5.7 -> blank
$ws.Cells.Item(12, 5).Value = -5.3          # E12: was blank -> -5.3
$ws.Cells.Item(14, 5).ClearContents()       # E14: was -5.4 -> blank
$ws.Cells.Item(20, 5).Value = -7.2          # E20: was blank -> -7.2
$ws.Cells.Item(21, 5).Value = -8.699999999999999  # E21: was blank -> -8.7
$ws.Cells.Item(22, 5).ClearContents()       # E22: was -6.1 -> blank
$ws.Cells.Item(23, 5).ClearContents()       # E23: was -7 -> blank

# --- Remove rows for samples "RM 232" (row 26) and "SC 92" (row 28) ---
# Delete bottom-up so row indices of earlier rows stay valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# --- After the row removal, old rows 33/34/35 (SC 132/SC 193/SC 232)
#     shift up to new rows 31/32/33; fix up their "missing data" cells
#     to match the new run ---
$ws.Cells.Item(31, 5).Value = -8.1          # E31 (was SC 132 row 33): blank -> -8.1
$ws.Cells.Item(32, 2).ClearContents()       # B32 (was SC 193 row 34): -19.9 -> blank
$ws.Cells.Item(33, 5).Value = -10.7         # E33 (was SC 232 row 35): blank -> -10.7
